# issue #5: add legislator_id, name, date into dataframe
#
# The 股票 (stocks) worksheet gets three new trailing columns appended to
# every row of its table: date, legislator_name, legislator_id.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$legislatorName = "黃志雄"
$legislatorId = 1366
$reportDate = "2011-12-21"
$lastRow = 13

# Header row (row 1) - match the existing header formatting (bold, centered,
# bordered) by copying it from an existing header cell before writing text.
$ws.Range("G1").Copy()
$ws.Range("H1:J1").PasteSpecial(-4122)

$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"

# Data rows 2-13 get the same three values repeated down the column.
# The date column is force-typed as text (leading apostrophe) so it is
# stored as the literal string "2011-12-21" rather than being auto-parsed
# into a date serial number.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = "'" + $reportDate
    $ws.Cells.Item($r, 9).Value = $legislatorName
    $ws.Cells.Item($r, 10).Value = $legislatorId
}
